$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    and push the existing 2022-Q3 summary row down to row 3.
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows("2:2").Insert()

# Give A2 the same numeric style ("s=2") the row had before the insert
# shifted it down to A3 (A3 still carries the original formatting).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A3").Value = 1

# ------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q3" detail sheet. The copy keeps the
#    old data (and becomes the new "2022-Q3" tab); the original sheet
#    object is turned into the new "2022-Q4" tab (this preserves the
#    rId2 / sheetId=2 slot for "2022-Q4" and hands the new rId/sheetId
#    to "2022-Q3", matching the target workbook layout).
# ------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("2022-Q3")
$oldQ3.Copy($null, $oldQ3)

$q4 = $oldQ3
$q3 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"
$q3.Name = "2022-Q3"

# ------------------------------------------------------------------
# 3. Turn $q4 into the 2022-Q4 detail sheet: drop the extra data rows
#    and fill in the two Q4 fund rows.
# ------------------------------------------------------------------
$q4.Rows("4:11").Delete()

# Re-style header row + index column to match the look used by
# freshly-authored sheets in this workbook (same style as the summary
# sheet's header / index column).
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q4.Range("A2:A3").PasteSpecial(-4122)

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'970073"
$q4.Range("C2").Value = "东证融汇成长优选混合A"
$q4.Range("D2").Value = "'0.38"
$q4.Range("E2").Value = "'89.59"
$q4.Range("F2").Value = "'0.82"
$q4.Range("G2").Value = "'0.0031"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'970074"
$q4.Range("C3").Value = "东证融汇成长优选混合C"
$q4.Range("D3").Value = "'0.11"
$q4.Range("E3").Value = "'89.59"
$q4.Range("F3").Value = "'0.82"
$q4.Range("G3").Value = "'0.0009"
$q4.Range("H3").Value = 8

# Drop the "typed as text via leading apostrophe" quote-prefix marker
# so these cells end up with no explicit style, matching plain
# inline-string cells elsewhere in the workbook.
$q4.Range("B2:G3").ClearFormats()
